$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Data")
$ws.Range("AW1").EntireColumn.Insert()
$ws.Range("AW6:AW8").Merge()
Write-Host "Merged AW6:AW8"
